# Applies two edits to "Week 1 Day 1.pptx":
#  1. Slide 3 title: split "Week 1 Learning Objectives" into
#     "Learning " + "Objectives" (dropping the leading "Week 1 "),
#     left-align the paragraph.
#  2. Slide 4 content placeholder: remove the "CU-Harvard referencing
#     guides" bullet paragraph entirely, and split the "Boolean algebra
#     and logic. " run into "Boolean " + "algebra and logic. ".

$p = $ppt.ActivePresentation

# --- Slide 3: Title -------------------------------------------------
$slide3 = $p.Slides.Item(3)
$title = $slide3.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange

$titleRange.Text = "Learning Objectives"
$titleRange.ParagraphFormat.Alignment = 1

$learningRun = $titleRange.Characters(1, 9)
$learningRun.Font.Bold = $true

$objectivesRun = $titleRange.Characters(10, 10)
$objectivesRun.Font.Bold = $true

# --- Slide 4: Content placeholder ------------------------------------
$slide4 = $p.Slides.Item(4)
$content = $slide4.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange

# Remove the "CU-Harvard referencing guides" paragraph (6th paragraph)
$cuHarvardPara = $contentRange.Paragraphs(6)
$cuHarvardPara.Delete()

# Split the "Boolean algebra and logic. " run (now the 6th paragraph)
$contentRange = $content.TextFrame.TextRange
$booleanPara = $contentRange.Paragraphs(6)
$booleanRun = $booleanPara.Characters(1, 8)
$booleanRun.Text = "Boolean "
